$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (28-Nov-2022) is inserted as row 6, pushing the
# existing rows 6 (7-Dec-2021) and 7 (18-Mar-2022) down to rows 7 and 8.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new record's data.
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C6").Value = "Arica y Parinacota"
$ws.Range("D6").Value = 44893
$ws.Range("E6").Value = 15
$ws.Range("F6").Value = 100112030
$ws.Range("G6").Value = "Poroto granado"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 3300
$ws.Range("K6").Value = 1200
$ws.Range("L6").Value = 1300
$ws.Range("M6").Value = 1261
$ws.Range("N6").Value = "`$/kilo"
$ws.Range("O6").Value = "Región de Arica y Parinacota"
$ws.Range("P6").Value = 1261
$ws.Range("Q6").Value = 1
$ws.Range("R6").Value = "Hortaliza"
